$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project_Charter")

# Problem description text (merged B7:C10) - replace long intro paragraph
# with the new, shorter non-conformities paragraph.
$ws.Range("B7").Value = "MEDISA company has presented several non-conformities with one of its main customers, highlighting incorrect measurements in its flaghship product. Therefore, the team will seek to reduce variations in the process and improve customer satisfaction."

# Benefits text (merged C12:C16) - update savings figure
$ws.Range("C12").Value = "A savings of `$545,592.96 MXN per year will be realized by reducing variations in the process. Comply with best practices and avoid product quality impacts. Increase customer satisfaction."

# Operators team list (H7:H9) - update names
$ws.Range("H7").Value = "Geardo, A."
$ws.Range("H8").Value = "Mario, M."
$ws.Range("H9").Value = "Citlali, B."

# Update the last active selection to match the saved state
$ws.Range("K15").Select()
